# Updated symbol list on Sun Dec 18 06:57:24 UTC 2022 with GitHub Actions
#
# Refreshes the cryptocurrency price/metadata sheet: most rows just get a
# new "Price" (column D) quote, while rows 41-43 also had their coin
# ordering reshuffled (Kick/BKEX/CEJI rotated) along with their Link/Coin
# and "Volume(1h)"-style label columns.
#
# NOTE: column D values are stored as literal text in this sheet (not real
# numbers), so a leading apostrophe is used when assigning them through
# COM to keep Excel from re-interpreting the numeric-looking strings as
# Number cells (which would also introduce floating point drift, e.g.
# 247.23 -> 247.22999999999999). Re-applying the "Normal" style right
# after keeps the cell format byte-for-byte the same as before the edit
# (Excel otherwise flags quote-prefixed cells with their own style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

# --- Column D "Price" refreshes -------------------------------------------------
Set-TextValue "D2"  "247.23"
Set-TextValue "D3"  "22.78"
Set-TextValue "D4"  "5.554"
Set-TextValue "D5"  "0.05616"
Set-TextValue "D8"  "0.8029"
Set-TextValue "D10" "0.1424"
Set-TextValue "D11" "0.07390"
Set-TextValue "D12" "0.03185"
Set-TextValue "D13" "0.02958"
Set-TextValue "D14" "0.09253"
Set-TextValue "D15" "0.001660"
Set-TextValue "D16" "3.094"
Set-TextValue "D17" "0.04713"
Set-TextValue "D18" "0.0005756"
Set-TextValue "D19" "0.006269"
Set-TextValue "D20" "0.001053"
Set-TextValue "D21" "0.003820"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "0.0004604"
Set-TextValue "D24" "3.981"
Set-TextValue "D25" "2.118"
Set-TextValue "D27" "0.1291"

# --- Row 41: now BKEXToken (was KickToken) --------------------------------------
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1045"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# --- Row 42: now CEJI (was BKEXToken) -------------------------------------------
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002972"
$ws.Range("E42").Value = "41CEJICEJI"

# --- Row 43: now KickToken (was CEJI) -------------------------------------------
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003249"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- Remaining Price / label refreshes ------------------------------------------
Set-TextValue "D44" "0.009823"
Set-TextValue "D45" "0.00005640"
Set-TextValue "D47" "0.6807"
Set-TextValue "D48" "0.02863"
$ws.Range("E48").Value = "47BOLOBOLO"
Set-TextValue "D49" "0.00002102"
